# Apply the "Modified Reg iExam TC's" edit:
#  - rename the sheet from "users" to "Worksheet"
#  - refresh a batch of generated credential strings
#  - update the numeric candidate id
#  - clear the bold/border formatting that used to decorate rows 1-2
#  - collapse the selection down to just A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Worksheet"

# Update the regenerated credential / id values
$ws.Range("A2").Value = "ZhyoF284"
$ws.Range("B2").Value = 23100602
$ws.Range("C2").Value = "qfzzgjw49"
$ws.Range("D2").Value = "A&3sP!h6"
$ws.Range("F2").Value = "xgWJZTBa"
$ws.Range("G2").Value = "ydGc"

# Strip the bold header / bordered formatting back to the default style
$ws.Range("A1:H2").ClearFormats()

# Collapse the selection to A1
$selResult = $ws.Range("A1").Select()
